# 自动更新Excel文件 - 每日库存/到期天数刷新脚本
# 规则：
#   对每一数据行 (第2行起)：
#     D = 总天数, E = 剩余天数, F = 开始时间 (yyyyMMdd)
#   若剩余天数 E 大于 1，则今天过去一天，剩余天数减 1，开始时间 F 不变。
#   若剩余天数 E 等于 1（即今天是最后一天），说明周期已满，自动续期：
#     剩余天数 E 重置为总天数 D，开始时间 F 更新为新的开始日期 20251014。
#   若某行的开始时间格式不是合法的 8 位日期（脏数据），则跳过该行，不做任何改动。

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newStartDate = 20251014

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 99 }

for ($r = 2; $r -le $lastRow; $r++) {
    $dVal = $ws.Range("D$r").Value2
    $eVal = $ws.Range("E$r").Value2
    $fVal = $ws.Range("F$r").Value2

    if ($null -eq $dVal -or $null -eq $eVal -or $null -eq $fVal) {
        continue
    }

    # 校验开始时间是否为合法的 8 位日期数字（例如 20251006），否则跳过该行
    $fText = [string][int64]$fVal
    if ($fText.Length -ne 8) {
        continue
    }

    $total = [int]$dVal
    $remaining = [int]$eVal

    if ($remaining -gt 1) {
        $ws.Range("E$r").Value = $remaining - 1
    } else {
        $ws.Range("E$r").Value = $total
        $ws.Range("F$r").Value = $newStartDate
    }
}
